$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the worker "MARLON ENRIQUE LEON ANTEQUERA" (row 16) from the
# statement. Deleting the row shifts every following row up by one
# (17->16, 18->17, 23->22, 24->23) and carries each row's own formatting
# along with it, which matches the structure of the updated workbook.
$ws.Rows("16").Delete()

# After the delete, the former "period 2111" row (old row 17) now sits at
# row 16 and the former "period 2110" row (old row 18) now sits at row 17.
# The refreshed statement lists period 2110 first, so swap the Periodo
# Mora / Valor Mora / Salario Basico values between the two remaining
# rows for YESICA PAOLA GUETO BARRIOS.
$periodo16 = $ws.Range("E16").Value2
$periodo17 = $ws.Range("E17").Value2
$valorMora16 = $ws.Range("F16").Value2
$valorMora17 = $ws.Range("F17").Value2
$salario16 = $ws.Range("G16").Value2
$salario17 = $ws.Range("G17").Value2

$ws.Range("E16").Value = $periodo17
$ws.Range("E17").Value = $periodo16
$ws.Range("F16").Value = $valorMora17
$ws.Range("F17").Value = $valorMora16
$ws.Range("G16").Value = $salario17
$ws.Range("G17").Value = $salario16

# Refresh the summary totals for the updated account statement.
$ws.Range("E11").Value = 109023
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2

# Column D ("Nombre Trabajador") was sized to fit the longest name; now
# that the longer "MARLON ENRIQUE LEON ANTEQUERA" row is gone, re-fit it
# to the remaining, shorter content.
$ws.Range("D1").ColumnWidth = 29
